$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# xlShiftDown = -4121, xlFormatFromLeftOrAbove = 0
$xlShiftDown = -4121
$xlFormatFromLeftOrAbove = 0

# The "list of origin countries" in column A is sorted alphabetically and
# ends at row 61 (Zimbabwe). Append two new countries ("few and early
# elections" decision samples) at rows 62-63, using Insert() so the new
# cells inherit the same formatting (style) as the row above them, just
# like the existing rows.
$ws.Range("A62").Insert($xlShiftDown, $xlFormatFromLeftOrAbove)
$ws.Range("A62").Value = "Comoros"

$ws.Range("A63").Insert($xlShiftDown, $xlFormatFromLeftOrAbove)
$ws.Range("A63").Value = "Kazakhstan"

# Mirror the author's view-state change: selection parked just past the
# new last row.
$ws.Range("A65").Select()

